$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.157.40'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.131.20'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.72'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.49'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.86%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.119.68'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.516'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.145'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.18%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.456'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.95%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.17'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.644.65'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.124.51'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.128.38'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.68'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '473.39'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.19'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.697'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.69'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '86.71'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.97'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.59%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.72'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.11'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.96'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.90%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.71'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.01%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.97%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.81'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.04'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0711'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '423.68'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.54%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.23'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.71'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -9.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.896.07'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.90%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.13'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.61'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.44%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.70%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.41'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.19%  '
